# Rename the original sheet and add a new "weights" sheet carrying the
# pollster MAE-based weighting table, per the commit:
# "Add methodology narrative and publish pollster MAE-based weights on dashboard"

$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> weighted_time_series (data/content unchanged) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "weighted_time_series"

# --- Add a new "weights" sheet right after it ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "weights"

# Header row (bold/centered header style matching the original sheet's headers)
$headers = @("조사기관", "mae", "weight", "weight_pct")
for ($c = 1; $c -le $headers.Count; $c++) {
    $ws2.Cells.Item(1, $c).Value = $headers[$c - 1]
}
$ws1.Range("A1:D1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)

# Pollster MAE / weight / weight_pct rows
$weightsData = @(
    @("리서치앤리서치", 1.715165533504263, 0.1852563467633035, 18.52563467633035),
    @("엠브레인퍼블릭", 1.988068392353487, 0.1598261418236138, 15.98261418236138),
    @("리서치뷰", 2.520935456058236, 0.1260426164691112, 12.60426164691112),
    @("에이스리서치", 2.706784148984008, 0.1173884888274516, 11.73884888274516),
    @("한국리서치", 2.984805325724369, 0.106454279645263, 10.6454279645263),
    @("조원씨앤아이", 3.005944774053744, 0.1057056349052044, 10.57056349052044),
    @("알앤써치", 4.15182333959994, 0.07653150792826108, 7.653150792826108),
    @("리얼미터", 4.376554982336593, 0.0726016929100001, 7.26016929100001),
    @("코리아리서치인터내셔널", 6.330433733753978, 0.05019329072779152, 5.019329072779152)
)

for ($r = 0; $r -lt $weightsData.Count; $r++) {
    $rowVals = $weightsData[$r]
    $rowIndex = $r + 2
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws2.Cells.Item($rowIndex, $c + 1).Value = $rowVals[$c]
    }
}

# Keep the originally-active sheet selected/active, as before.
$ws1.Activate()
